$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Calculate confidence intervals when the population variance is unknown:
# clear the sample-mean and standard-error helper formulas in E9/E11 so the
# sheet now reflects the "unknown variance" scenario (dependent CI bounds
# recalculate to 0 until new inputs are supplied).
$ws.Range("E9").ClearContents()
$ws.Range("E11").ClearContents()

# Widen column E to make room for the (now empty) inputs.
$ws.Columns.Item(5).ColumnWidth = 60

# Bump the row heights for rows 9 and 11 slightly.
$ws.Rows.Item(9).RowHeight = 12.8
$ws.Rows.Item(11).RowHeight = 12.8

# Move the active selection to E10.
[void]$ws.Range("E10").Select()
